# Refresh of the "cryptos" price/volume table (GitHub Actions scheduled
# update). Updates Price (column D) and Volume(1h) (column E) for most
# rows, and additionally swaps the BinanceUSD / EthereumClassic rows
# (28 <-> 29) along with their own refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be forced back to
# Text storage (the source data stores prices/volumes as text, e.g. "212.27"),
# otherwise Excel's auto-detection would silently store them as numeric values.
function Set-TextValue($cell, $text) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "27.640.42"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.634.13"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue "D5" "212.27"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("E7").Value = "  +0.05%  "
Set-TextValue "D8" "23.32"
$ws.Range("E8").Value = "  +1.16%  "
Set-TextValue "D9" "0.265"
$ws.Range("E9").Value = "  +2.85%  "
$ws.Range("E10").Value = "  +0.26%  "
Set-TextValue "D11" "0.0858"
$ws.Range("E11").Value = "  -4.06%  "
$ws.Range("D12").Value = "1.864.95"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "1.639.41"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  +0.13%  "
Set-TextValue "D15" "0.553"
$ws.Range("E15").Value = "  -1.02%  "
Set-TextValue "D16" "65.23"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").Value = "27.618.62"
$ws.Range("E17").Value = "  -0.13%  "
Set-TextValue "D18" "230.69"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("E19").Value = "  -0.28%  "
Set-TextValue "D20" "7.58"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("E21").Value = "  +0.06%  "
Set-TextValue "D22" "10.62"
$ws.Range("E22").Value = "  +4.38%  "
$ws.Range("E23").Value = "  +1.39%  "
Set-TextValue "D24" "2.11"
$ws.Range("E24").Value = "  +4.38%  "
Set-TextValue "D25" "149.02"
$ws.Range("E25").Value = "  -1.05%  "
Set-TextValue "D26" "6.88"
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D28" "1.00"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D29" "15.51"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("D33").Value = "1.477.65"
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("E36").Value = "  -1.33%  "
Set-TextValue "D37" "0.959"
$ws.Range("E37").Value = "  +6.43%  "
$ws.Range("E38").Value = "  -0.06%  "
Set-TextValue "D39" "0.558"
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("E42").Value = "  +0.02%  "
Set-TextValue "D43" "67.79"
$ws.Range("E43").Value = "  -3.27%  "
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("E45").Value = "  -1.35%  "
Set-TextValue "D46" "5.33"
$ws.Range("E46").Value = "  -5.10%  "
$ws.Range("D47").Value = "1.774.77"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  +0.31%  "
Set-TextValue "D49" "87.80"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("E50").Value = "  -1.64%  "
Set-TextValue "D51" "0.0991"
$ws.Range("E51").Value = "  +0.02%  "
